$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "ser" numbers of the blog entries in the dashboard's last row
# (each entry bumped up by one, making room for the newest post, #121,
# "good things in Pakistan"):
#   I8: ser 118 -> 119
#   E8: ser 119 -> 120
#   C8: ser 120 -> 121
$ws.Range("I8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 119"
$ws.Range("E8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 120"
$ws.Range("C8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 121"
